# Message_List_20180321.xlsx - add three new request/response message pairs
# (add_transaction, get_content_list_by_user, get_content_list_owner, get_content)
# into the previously-empty rows 35-42 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 42 (the "get_content" response) gains a full row of fields, so first
# copy the formatting used by the other "response with content fields" rows
# (e.g. row 20 / row 28) onto row 42's F:P range before filling in values.
$ws.Range("F20:L20").Copy() | Out-Null
$ws.Range("F42:L42").PasteSpecial(-4122) | Out-Null
$ws.Range("M20:O20").Copy() | Out-Null
$ws.Range("M42:O42").PasteSpecial(-4122) | Out-Null
$ws.Range("P20").Copy() | Out-Null
$ws.Range("P42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 35: add_transaction (request: Web Client -> Server) ---
$ws.Range("B35").Value = "add_transaction"
$ws.Range("C35").Value = "Web Client"
$ws.Range("D35").Value = "Server"
$ws.Range("E35").Value = "Session Code"
$ws.Range("F35").Value = "User ID"
$ws.Range("G35").Value = "Amount"
$ws.Range("H35").Value = "Description"

# --- Row 36: add_transaction (response: Server -> Web Client) ---
$ws.Range("B36").Value = "add_transaction"
$ws.Range("C36").Value = "Server"
$ws.Range("D36").Value = "Web Client"
$ws.Range("E36").Value = "Message status"

# --- Row 37: get_content_list_by_user (request: Phone -> Server) ---
$ws.Range("B37").Value = "get_content_list_by_user"
$ws.Range("C37").Value = "Phone"
$ws.Range("D37").Value = "Server"
$ws.Range("E37").Value = "User ID"

# --- Row 38: get_content_list_by_user (response: Server -> Phone) ---
$ws.Range("B38").Value = "get_content_list_by_user"
$ws.Range("C38").Value = "Server"
$ws.Range("D38").Value = "Phone"
$ws.Range("E38").Value = "Message status"
$ws.Range("F38").Value = "Content list"

# --- Row 39: get_content_list_owner (request: Phone -> Server) ---
$ws.Range("B39").Value = "get_content_list_owner"
$ws.Range("C39").Value = "Phone"
$ws.Range("D39").Value = "Server"
$ws.Range("E39").Value = "Session Code"

# --- Row 40: get_content_list_owner (response: Server -> Phone) ---
$ws.Range("B40").Value = "get_content_list_owner"
$ws.Range("C40").Value = "Server"
$ws.Range("D40").Value = "Phone"
$ws.Range("E40").Value = "Message status"
$ws.Range("F40").Value = "Content list"

# --- Row 41: get_content (request: Phone -> Server) ---
$ws.Range("B41").Value = "get_content"
$ws.Range("C41").Value = "Phone"
$ws.Range("D41").Value = "Server"
$ws.Range("E41").Value = "Session Code (?)"
$ws.Range("F41").Value = "Content ID"

# --- Row 42: get_content (response: Server -> Phone) ---
$ws.Range("B42").Value = "get_content"
$ws.Range("C42").Value = "Server"
$ws.Range("D42").Value = "Phone"
$ws.Range("E42").Value = "Message status"
$ws.Range("F42").Value = "Title"
$ws.Range("G42").Value = "Content"
$ws.Range("H42").Value = "Address"
$ws.Range("I42").Value = "Stretch"
$ws.Range("J42").Value = "Price"
$ws.Range("K42").Value = "Priority"
$ws.Range("L42").Value = "Created date"
$ws.Range("M42").Value = "Expire date"
$ws.Range("N42").Value = "Image list"
$ws.Range("O42").Value = "Avatar"
$ws.Range("P42").Value = "Expand Data"

# --- Column H widened slightly now that "Description" is the longest entry ---
$ws.Columns("H").ColumnWidth = 10.36328125

# --- Selection / scroll position as saved in the final workbook ---
$ws.Range("A1").Select()
$ws.Range("B43").Select()
